$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 8: "Bolded Row" demo -----------------------------------------
# A8 carries the text and the bold style that the (unsupported-here)
# row-level default would have applied.
$ws.Range("A8").Value = "Bolded Row"
$ws.Range("A8").Font.Bold = $true

# C8 is an otherwise-empty cell that shows what the row(bold) + column
# C (italic) combination looks like -> bold-italic, no value.
$ws.Range("C8").Font.Bold = $true
$ws.Range("C8").Font.Italic = $true

# --- Row 9: "Bolded Column" demo ---------------------------------------
# B9 shows what column B's default (bold) formatting looks like on an
# otherwise normal row.
$ws.Range("B9").Value = "Bolded Column"
$ws.Range("B9").Font.Bold = $true

# --- Row 10: "Bolded Row Italic Column" demo ---------------------------
# C10 combines row bold + column italic explicitly (bold-italic).
$ws.Range("C10").Value = "Bolded Row Italic Column"
$ws.Range("C10").Font.Bold = $true
$ws.Range("C10").Font.Italic = $true

# --- Row 11: plain "Normal" cell in column D ---------------------------
$ws.Range("D11").Value = "Normal"

# --- Column widths/defaults --------------------------------------------
# Column B default = bold, column C default = italic, column D = plain.
$ws.Columns.Item(2).Font.Bold = $true
$ws.Columns.Item(3).Font.Italic = $true

$ws.Columns.Item(2).AutoFit() | Out-Null
$ws.Columns.Item(3).AutoFit() | Out-Null
$ws.Columns.Item(4).AutoFit() | Out-Null

# --- Selection matches the authored workbook ----------------------------
$ws.Range("C9").Select() | Out-Null
